$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - new globalization columns D, E, F (C3 unchanged)
$ws.Range("D3").Value = -1823.9272228992447
$ws.Range("E3").Value = -1721.1270833989984
$ws.Range("F3").Value = -1999.8353116116334

# Row 4 - C4 value updated, plus new D, E, F
$ws.Range("C4").Value = 5711.4209073743714
$ws.Range("D4").Value = 4650.488621732502
$ws.Range("E4").Value = 5279.5877045758834
$ws.Range("F4").Value = 5280.2566044854921

# Row 5 - new D, E, F (C5 unchanged)
$ws.Range("D5").Value = -23404.131545167438
$ws.Range("E5").Value = -1271.1217576267572
$ws.Range("F5").Value = -23404.131545167438

# Row 6 - new D, E, F (C6 unchanged)
$ws.Range("D6").Value = -75551.943901742939
$ws.Range("E6").Value = -75551.943901742939
$ws.Range("F6").Value = -75551.943901742939

# Row 7 - new D, E, F (C7 unchanged)
$ws.Range("D7").Value = -20987.965400176996
$ws.Range("E7").Value = 22886.218720182835
$ws.Range("F7").Value = -16793.762565254285

# Update the active selection to match the saved workbook state
$ws.Range("I10").Select()
